$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.698.50"
$ws.Range("E2").Value = "  +5.65%  "
$ws.Range("D3").Value = "1.729.48"
$ws.Range("E3").Value = "  +4.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'227.08"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'0.5434"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.2732"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "'0.06667"
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("D10").Value = "'21.79"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").Value = "'0.07777"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'4.681"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.721.62"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("D14").Value = "1.968.34"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").Value = "'0.5945"
$ws.Range("E15").Value = "  +5.22%  "
$ws.Range("D16").Value = "0.0₅8381"
$ws.Range("D17").Value = "'68.94"
$ws.Range("E17").Value = "  +4.67%  "
$ws.Range("D18").Value = "27.698.92"
$ws.Range("E18").Value = "  +5.68%  "
$ws.Range("E19").Value = "  +17.09%  "
$ws.Range("D20").Value = "'4.798"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D22").Value = "'10.84"
$ws.Range("E22").Value = "  +4.26%  "
$ws.Range("D23").Value = "'6.195"
$ws.Range("E23").Value = "  +3.09%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "'1.731"
$ws.Range("E26").Value = "  +13.39%  "
$ws.Range("D27").Value = "'0.1247"
$ws.Range("D28").Value = "'7.445"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +5.79%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").Value = "'3.658"
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("D33").Value = "'3.495"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").Value = "'1.670"
$ws.Range("E34").Value = "  +5.38%  "
$ws.Range("D35").Value = "'0.9724"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("D36").Value = "'2.843"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "'2.437"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'0.5969"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").Value = "'0.01664"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("D40").Value = "'5.909"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'0.8594"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "1.046.80"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'101.33"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "1.873.22"
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("E46").Value = "  +8.76%  "
$ws.Range("D47").Value = "'59.47"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").Value = "'8.242"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'0.4429"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "'0.05328"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'0.9997"
$ws.Range("E51").Value = "  -0.66%  "
